# Add a team record (Wins/Losses/Ties) to the roster sheet.
# New columns AD:AF are appended after the existing Salary/Unnamed:28 columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column titles, styled like the rest of the header row.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold font, borders, centered alignment) from an
# existing header cell onto the three new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Every player row (2 through 48) gets the team's season record.
for ($r = 2; $r -le 48; $r++) {
    $ws.Range("AD$r").Value = 88
    $ws.Range("AE$r").Value = 74
    $ws.Range("AF$r").Value = 0
}
